# Update the COVID-19 Valais figures sheet with the latest batch of daily
# data (rows 513-524): correct a few "new cases" counts that had been
# entered too low, fill in the running totals for 22.07-02.08.2021 that
# were still blank, and move the frozen-pane scroll/selection down to the
# newly-entered rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to previously-entered rows ---------------------------
# Column C holds the day's new cases; column B is a running total formula
# that recalculates automatically once C changes.
$ws.Range("C513").Value = 33
$ws.Range("C518").Value = 23
$ws.Range("C519").Value = 20
$ws.Range("C520").Value = 22

# --- Newly-entered rows -------------------------------------------------
# Row 521 - 30.07.2021
$ws.Range("C521").Value = 19
$ws.Range("E521").Value = 1
$ws.Range("F521").Value = 1
$ws.Range("G521").Value = 3
$ws.Range("L521").Value = 0
$ws.Range("M521").Value = 0

# Row 522 - 31.07.2021
$ws.Range("C522").Value = 17
$ws.Range("E522").Value = 1
$ws.Range("F522").Value = 1
$ws.Range("G522").Value = 2
$ws.Range("L522").Value = 0
$ws.Range("M522").Value = 0

# Row 523 - 01.08.2021
$ws.Range("C523").Value = 4
$ws.Range("E523").Value = 1
$ws.Range("F523").Value = 1
$ws.Range("G523").Value = 3
$ws.Range("L523").Value = 0
$ws.Range("M523").Value = 0

# Row 524 - 02.08.2021 (new-cases count for this day not filled in yet)
$ws.Range("E524").Value = 1
$ws.Range("F524").Value = 1
$ws.Range("G524").Value = 3
$ws.Range("L524").Value = 0
$ws.Range("M524").Value = 0

# --- Scroll the frozen pane down to the rows just edited ---------------
$null = $excel.Goto($ws.Range("B494"), $true)
$null = $ws.Range("O501").Select()
